$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "'0420172010219   "
$ws.Range("C3").Value = "'55309     "
$ws.Range("C2").Value = "'32610 "
$ws.Range("B2").Value = "'0420194406901 "

$ws.Range("D3").Select()
